$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update marking value (Right answer marks) from 3 to 5
$ws.Range("B11").Value = 5

# Update total marks obtained (Right * Marking) from 21 to 35
$ws.Range("B12").Value = 35

# Update correct/total marks label from "19/84" to "35/140"
$ws.Range("E12").Value = "35/140"
